# Add filter: a new per-day record sheet "2024-07-14" (same layout as the
# other daily record sheets) and append its summary row to the "current"
# sheet so the doctor filter can pick it up.

$wb = $excel.ActiveWorkbook

# --- 1. New worksheet "2024-07-14" with the standard record header row ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2024-07-14"

$headers = @("ID", "Время", "ФИО пациента", "Врач", "Врач_Индекс", "М\Ж\Р", "Дата рождения", "Причина", "Давление")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Move it to the end of the tab strip, after the last existing daily sheet.
$lastSheet = $wb.Worksheets.Item("2024-07-13")
$newSheet.Move($null, $lastSheet)

# Match the outline/page-setup boilerplate used by the other daily sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$ps = $newSheet.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# --- 2. Append the new day's summary row to the "current" sheet ---
$current = $wb.Worksheets.Item("current")

$dateCell = $current.Range("A8")
# Force the date-looking label to be stored as plain text (not an Excel
# serial date number), then drop the temporary Text number format so the
# cell keeps the sheet's default (unstyled) look, matching the other rows.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-07-14"
$dateCell.ClearFormats()

$current.Range("B8").Value = -1
